$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Unprotect()

# Update the "as of" date in the confidential disclaimer text (A10)
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-27 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-7
$ws.Range("D2").Value = 0.477583239259485
$ws.Range("E2").Value = 0.001951600312255941

$ws.Range("D3").Value = 0.3392667819876123
$ws.Range("E3").Value = 0.002443609022556359

$ws.Range("D4").Value = 0.09779674039480971
$ws.Range("E4").Value = -0.003686797752809112

$ws.Range("D5").Value = 0.05350098815865235
$ws.Range("E5").Value = -0.0004584527220631296

$ws.Range("D6").Value = 0.03185225019944065
$ws.Range("E6").Value = -0.002598902685532867

$ws.Range("D7").Value = 0.9999999999999999
$ws.Range("E7").Value = 0.00129322159342804

# Restore sheet protection (it was unprotected above only to allow the edits)
$ws.Protect()
